$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 13: B13=1, C13=2, D13=formula B13+C13
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 2
$ws.Range("D13").Formula = "=B13+C13"

# Conditional formatting on D13: highlight green ("Good" style) when value > 0
$fc = $ws.Range("D13").FormatConditions.Add(1, 5, "0")
$fc.Font.Color = 24832
$fc.Interior.Color = 13561798

# Show formulas in the sheet view, and select D13
$excel.ActiveWindow.DisplayFormulas = $true
$ws.Range("D13").Select() | Out-Null

# Adjust column widths (A is new, B/C/D narrower, no longer "best fit")
$ws.Columns.Item(1).ColumnWidth = 5.451822916666667
$ws.Columns.Item(2).ColumnWidth = 20.736979166666668
$ws.Columns.Item(3).ColumnWidth = 12.166666666666666
$ws.Columns.Item(4).ColumnWidth = 10.592447916666666
